$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8").Value = 45995
$ws.Range("D8").Value = 163.51
$ws.Range("E8").Value = 161.86
$ws.Range("F8").Value = 171.86
$ws.Range("G8").Value = 161.97

$ws.Range("A9").Value = 45995
$ws.Range("D9").Value = 163.51
$ws.Range("E9").Value = 161.86
$ws.Range("F9").Value = 171.86
$ws.Range("G9").Value = 161.97

$ws.Range("A10").Value = 45995
$ws.Range("D10").Value = 165.38
$ws.Range("E10").Value = 164.6
$ws.Range("F10").Value = 174.6
$ws.Range("G10").Value = 165.12

$ws.Range("A11").Value = 45994
$ws.Range("D11").Value = 163.91
$ws.Range("E11").Value = 161.2
$ws.Range("F11").Value = 171.2
$ws.Range("G11").Value = 161.31

$ws.Range("A12").Value = 45994
$ws.Range("D12").Value = 163.91
$ws.Range("E12").Value = 161.2
$ws.Range("F12").Value = 171.2
$ws.Range("G12").Value = 161.31

$ws.Range("A13").Value = 45994
$ws.Range("D13").Value = 165.71
$ws.Range("E13").Value = 164
$ws.Range("F13").Value = 174
$ws.Range("G13").Value = 164.51

$ws.Range("A17").Value = 45995
$ws.Range("D17").Value = 168.37
$ws.Range("E17").Value = 167.07
$ws.Range("F17").Value = 177.07

$ws.Range("A18").Value = 45994
$ws.Range("D18").Value = 168.77
$ws.Range("E18").Value = 166.31
$ws.Range("F18").Value = 176.31

$ws.Range("A22").Value = 45995
$ws.Range("D22").Value = 164.32
$ws.Range("E22").Value = 163.79
$ws.Range("F22").Value = 173.39
$ws.Range("G22").Value = 164.95

$ws.Range("A23").Value = 45995
$ws.Range("D23").Value = 170.4
$ws.Range("E23").Value = 168.83
$ws.Range("F23").Value = 178.83

$ws.Range("A24").Value = 45995
$ws.Range("D24").Value = 170.16
$ws.Range("E24").Value = 169.12
$ws.Range("F24").Value = 179.12

$ws.Range("A25").Value = 45995
$ws.Range("D25").Value = 170.77
$ws.Range("E25").Value = 168.55
$ws.Range("F25").Value = 178.55
$ws.Range("G25").Value = 168.32

$ws.Range("A26").Value = 45995
$ws.Range("D26").Value = 169.59
$ws.Range("E26").Value = 170.01
$ws.Range("F26").Value = 180.01

$ws.Range("A27").Value = 45994
$ws.Range("D27").Value = 164.94
$ws.Range("E27").Value = 163.21
$ws.Range("F27").Value = 172.81
$ws.Range("G27").Value = 164.37

$ws.Range("A28").Value = 45994
$ws.Range("D28").Value = 170.5
$ws.Range("E28").Value = 168.01
$ws.Range("F28").Value = 178.01

$ws.Range("A29").Value = 45994
$ws.Range("D29").Value = 170.27
$ws.Range("E29").Value = 168.49
$ws.Range("F29").Value = 178.49

$ws.Range("A30").Value = 45994
$ws.Range("D30").Value = 171.1
$ws.Range("E30").Value = 167.91
$ws.Range("F30").Value = 177.91
$ws.Range("G30").Value = 167.68

$ws.Range("A31").Value = 45994
$ws.Range("D31").Value = 169.72
$ws.Range("E31").Value = 169.37
$ws.Range("F31").Value = 179.37

$ws.Range("A35").Value = 45995
$ws.Range("D35").Value = 163.75
$ws.Range("E35").Value = 161.6
$ws.Range("F35").Value = 170.6

$ws.Range("A36").Value = 45994
$ws.Range("D36").Value = 164.07
$ws.Range("E36").Value = 160.99
$ws.Range("F36").Value = 169.99

$ws.Range("A40").Value = 45995
$ws.Range("D40").Value = 169.62
$ws.Range("E40").Value = 167.59
$ws.Range("F40").Value = 177.59

$ws.Range("A41").Value = 45995
$ws.Range("D41").Value = 169.33
$ws.Range("E41").Value = 168.01
$ws.Range("F41").Value = 178.01

$ws.Range("A42").Value = 45994
$ws.Range("D42").Value = 169.89
$ws.Range("E42").Value = 166.7
$ws.Range("F42").Value = 176.7

$ws.Range("A43").Value = 45994
$ws.Range("D43").Value = 169.6
$ws.Range("E43").Value = 167.12
$ws.Range("F43").Value = 177.12

$ws.Range("A47").Value = 45995
$ws.Range("D47").Value = 164.33
$ws.Range("E47").Value = 162.62
$ws.Range("F47").Value = 172.62

$ws.Range("A48").Value = 45995
$ws.Range("D48").Value = 164.19
$ws.Range("E48").Value = 162.73
$ws.Range("F48").Value = 172.73

$ws.Range("A49").Value = 45994
$ws.Range("D49").Value = 164.4
$ws.Range("E49").Value = 162.04
$ws.Range("F49").Value = 172.04

$ws.Range("A50").Value = 45994
$ws.Range("D50").Value = 164.29
$ws.Range("E50").Value = 162.16
$ws.Range("F50").Value = 172.16

$ws.Range("A54").Value = 45995
$ws.Range("D54").Value = 179.27
$ws.Range("E54").Value = 179.13
$ws.Range("F54").Value = 189.13

$ws.Range("A55").Value = 45995
$ws.Range("D55").Value = 167.35
$ws.Range("E55").Value = 174.41
$ws.Range("F55").Value = 184.41

$ws.Range("A56").Value = 45995
$ws.Range("D56").Value = 169.7

$ws.Range("A57").Value = 45995
$ws.Range("D57").Value = 168.75
$ws.Range("E57").Value = 168.67

$ws.Range("A58").Value = 45995
$ws.Range("D58").Value = 164.66
$ws.Range("E58").Value = 164.73
$ws.Range("F58").Value = 174.73

$ws.Range("A59").Value = 45995
$ws.Range("D59").Value = 171.27
$ws.Range("E59").Value = 176.52

$ws.Range("A60").Value = 45994
$ws.Range("D60").Value = 179.65
$ws.Range("E60").Value = 178.22
$ws.Range("F60").Value = 188.22

$ws.Range("A61").Value = 45994
$ws.Range("D61").Value = 167.68
$ws.Range("E61").Value = 173.88
$ws.Range("F61").Value = 183.88

$ws.Range("A62").Value = 45994
$ws.Range("D62").Value = 170.01

$ws.Range("A63").Value = 45994
$ws.Range("D63").Value = 169.17
$ws.Range("E63").Value = 168.15

$ws.Range("A64").Value = 45994
$ws.Range("D64").Value = 165.07
$ws.Range("E64").Value = 164.2
$ws.Range("F64").Value = 174.2

$ws.Range("A65").Value = 45994
$ws.Range("D65").Value = 171.69
$ws.Range("E65").Value = 175.75
